$d = $word.ActiveDocument

# Locate the paragraph containing the last bibliography entry
# ("1999. Loeb, A. Biological monitoring ..."); this paragraph must be kept.
$bibIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*1999. Loeb, A. Biological monitoring of aquatic systems. McGraw-Hill. 1998.*") {
        $bibIndex = $i
    }
}

# Locate the paragraph containing the site footer/copyright text
# ("... Powered by Jekyll and Github pages ..."); this paragraph (and
# everything between it and the bibliography paragraph above) must be removed.
$footerIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Powered by Jekyll and Github pages*") {
        $footerIndex = $i
    }
}

if ($bibIndex -ne $null -and $footerIndex -ne $null) {
    $startPara = $d.Paragraphs.Item($bibIndex + 1)
    $endPara = $d.Paragraphs.Item($footerIndex)
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
